$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# Update the project website URL (underscore instead of hyphen)
$ws.Range("B4").Value = "https://github.com/RaphaelNajera/Sunlight_Sensor"

# Activate the sheet and update the selection/view to B4
$ws.Activate()
$ws.Range("B4").Select()
